$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9; this shifts existing rows 9-53 down to 10-54,
# carrying all their values/styles along automatically.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44575
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112012
$ws.Range("G9").Value = "Espinaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 11000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11600
$ws.Range("N9").Value = "$/cuna 10 kilos"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 1160
$ws.Range("Q9").Value = 10
$ws.Range("R9").Value = "Hortaliza"
